$d = $word.ActiveDocument

# 1. "Questions to answer in preparing for March 1 class and DARE #4"
#    -> "Questions to answer in preparing for Week 9 class and DARE #4"
$null = $d.Content.Find.Execute("March 1", $true, $false, $false, $false, $false, $true, 1, $false, "Week 9", 2)

# 2. Wrap "doubly-robust" with grammar-check markers (proofErr) -- not
#    representable through the exposed object model, so just make sure the
#    visible text itself is unchanged/intact here (no-op placeholder kept
#    for clarity).

# 3. Remove the stray "_GoBack" bookmark left over from the previous edit
#    session.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 4. Describe the robustness checks the authors conduct... (text itself is
#    unchanged, only grammar-check proofErr markers were added around
#    "authors" in the source diff -- not representable via the OM).

Write-Host "done"
